$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 22.59487733333333
$ws.Cells.Item(2, 8).Value = 67.784632
$ws.Cells.Item(2, 9).Value = 0.7395019553569895
$ws.Cells.Item(2, 10).Value = 0.7395019553569895
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 0.1030763333333333
$ws.Cells.Item(2, 14).Value = 0.309229
$ws.Cells.Item(2, 15).Value = 0.01126512502660735
$ws.Cells.Item(2, 16).Value = 0.01126512502660735
$ws.Cells.Item(2, 17).Value = 2.328997107636444
$ws.Cells.Item(2, 18).Value = 20.960973968728
$ws.Cells.Item(2, 19).Value = 0.008330581984517095
$ws.Cells.Item(2, 20).Value = 0.008330581984517096
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 22.59487733333333
$ws.Cells.Item(3, 8).Value = 67.784632
$ws.Cells.Item(3, 9).Value = 0.7395019553569895
$ws.Cells.Item(3, 10).Value = 0.7395019553569895
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 9.046962666666667
$ws.Cells.Item(3, 14).Value = 27.140888
$ws.Cells.Item(3, 15).Value = 0.9887348749733926
$ws.Cells.Item(3, 16).Value = 0.9887348749733927
$ws.Cells.Item(3, 17).Value = 204.4150116925796
$ws.Cells.Item(3, 18).Value = 1839.735105233216
$ws.Cells.Item(3, 19).Value = 0.7311713733724724
$ws.Cells.Item(3, 20).Value = 0.7311713733724725
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 0.3045986666666667
$ws.Cells.Item(4, 8).Value = 0.913796
$ws.Cells.Item(4, 9).Value = 0.00996913177602551
$ws.Cells.Item(4, 10).Value = 0.00996913177602551
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.1030763333333333
$ws.Cells.Item(4, 14).Value = 0.309229
$ws.Cells.Item(4, 15).Value = 0.01126512502660735
$ws.Cells.Item(4, 16).Value = 0.01126512502660735
$ws.Cells.Item(4, 17).Value = 0.03139691369822222
$ws.Cells.Item(4, 18).Value = 0.282572223284
$ws.Cells.Item(4, 19).Value = 0.0001123035158636516
$ws.Cells.Item(4, 20).Value = 0.0001123035158636516
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 0.3045986666666667
$ws.Cells.Item(5, 8).Value = 0.913796
$ws.Cells.Item(5, 9).Value = 0.00996913177602551
$ws.Cells.Item(5, 10).Value = 0.00996913177602551
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 9.046962666666667
$ws.Cells.Item(5, 14).Value = 27.140888
$ws.Cells.Item(5, 15).Value = 0.9887348749733926
$ws.Cells.Item(5, 16).Value = 0.9887348749733927
$ws.Cells.Item(5, 17).Value = 2.755692765649778
$ws.Cells.Item(5, 18).Value = 24.801234890848
$ws.Cells.Item(5, 19).Value = 0.009856828260161859
$ws.Cells.Item(5, 20).Value = 0.009856828260161859
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 7.654706
$ws.Cells.Item(6, 8).Value = 22.964118
$ws.Cells.Item(6, 9).Value = 0.2505289128669849
$ws.Cells.Item(6, 10).Value = 0.2505289128669849
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.1030763333333333
$ws.Cells.Item(6, 14).Value = 0.309229
$ws.Cells.Item(6, 15).Value = 0.01126512502660735
$ws.Cells.Item(6, 16).Value = 0.01126512502660735
$ws.Cells.Item(6, 17).Value = 0.7890190272246667
$ws.Cells.Item(6, 18).Value = 7.101171245021999
$ws.Cells.Item(6, 19).Value = 0.002822239526226605
$ws.Cells.Item(6, 20).Value = 0.002822239526226605
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 7.654706
$ws.Cells.Item(7, 8).Value = 22.964118
$ws.Cells.Item(7, 9).Value = 0.2505289128669849
$ws.Cells.Item(7, 10).Value = 0.2505289128669849
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 9.046962666666667
$ws.Cells.Item(7, 14).Value = 27.140888
$ws.Cells.Item(7, 15).Value = 0.9887348749733926
$ws.Cells.Item(7, 16).Value = 0.9887348749733927
$ws.Cells.Item(7, 17).Value = 69.25183940630934
$ws.Cells.Item(7, 18).Value = 623.266554656784
$ws.Cells.Item(7, 19).Value = 0.2477066733407583
$ws.Cells.Item(7, 20).Value = 0.2477066733407583
